# Saldo_guide.xlsx update
# - Rename the single worksheet to reflect the new export timestamp
# - Bump every "date" value in column G (rows 2-274) from 45617 to 45618
# - Correct two balance rows (51 and 120) whose Amount (E) / Balance (H)
#   values were re-extracted with new figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to the new export id
$ws.Name = "IClientBalance-20241122-105326-"

# Find the last used row/column (dimension is A1:H274)
$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

# Column G holds the snapshot date for every data row (2..274); shift it by one day
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 45618
}

# Row 51: Amount / Balance corrected from 25411.56 to 93.38
$ws.Cells.Item(51, 5).Value = 93.38
$ws.Cells.Item(51, 8).Value = 93.38

# Row 120: Amount / Balance corrected from 897.44 to 22823.18
$ws.Cells.Item(120, 5).Value = 22823.18
$ws.Cells.Item(120, 8).Value = 22823.18
